$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 9 ----
$ws.Cells.Item(9,1).Value2  = "np_harvest_individuals_arc2016"
$ws.Cells.Item(9,2).Value2  = "No. individuals captured by country (seals for pelts, walrus and narwhal for ivory)"
$ws.Cells.Item(9,3).Value2  = "FAO marine mammal capture data"
$ws.Cells.Item(9,7).Value2  = "Estimated total mortality of individuals per year"
$ws.Cells.Item(9,8).Value2  = "http://www.nmfs.noaa.gov/pr/sars/pdf/alaska2014_summary_final.pdf"
$ws.Cells.Item(9,9).Value2  = "FAO catch statistics. Large seal population around Jan Mayen but catch reporting goes to Norway. Seal stocks are managed jointly by Russia/Norway and quotas from ICES."
$ws.Cells.Item(9,11).Value2 = "V small licenced catch in Svalbard.  Can't get statistics so N/A"
$ws.Cells.Item(9,12).Value2 = "N/A"
$ws.Cells.Item(9,13).Value2 = "FAO statistics. "

# ---- Row 10 ----
$ws.Cells.Item(10,1).Value2  = "NP_quotas"
$ws.Cells.Item(10,2).Value2  = "Quotas set for NP products"
$ws.Cells.Item(10,3).Value2  = "Quotas set jointly by Russia/Norway as advised by ICES- still trying to get hold"
$ws.Cells.Item(10,7).Value2  = "No quotas set - have limited Potential Biological Removal data for a couple of species"
$ws.Cells.Item(10,9).Value2  = "Quotas set jointly by Russia/Norway as advised by ICES- still trying to get hold"
$ws.Cells.Item(10,11).Value2 = "N/A"
$ws.Cells.Item(10,13).Value2 = "No quotas for seal - obtained quotas for Walrus and 2015/2016 for narwhal - trying to get further back in time"

# ---- Hyperlink on H9 ----
$h9 = $ws.Cells.Item(9,8)
$h9.Hyperlinks.Add($h9, "http://www.nmfs.noaa.gov/pr/sars/pdf/alaska2014_summary_final.pdf") | Out-Null
$h9.Style = "Hyperlink"
$h9.WrapText = $true

# ---- Row heights ----
$ws.Rows.Item(9).RowHeight = 158.4
$ws.Rows.Item(10).RowHeight = 115.2

# ---- View state ----
$ws.Application.ActiveWindow.ScrollRow = 9
$ws.Range("O10").Select() | Out-Null
